$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column U (shifting nothing, since it's past the last used column)
# using "copy formatting from left" so the new cells inherit the same style that
# column T already uses (style index referenced by T2:T14 / S2:S14).
$null = $ws.Range("U1:U14").Insert(-4161, 0)

# Header for the new column.
$ws.Range("U1").Value = "QOIList"

# Rows 2-4 (Bernoulli-Pi, Bernoulli-Logit, Bernoulli-Logit-X) get the 2-item QOI list.
$ws.Range("U2:U4").Value = 'list("Predicted Values", "Expected Values")'

# Rows 5-14 (remaining distributions) get the 3-item QOI list.
$ws.Range("U5:U14").Value = 'list("Predicted Values", "Expected Values", "Probability Y > 1")'

# Update the active selection / view to the new column's header cell.
$null = $ws.Range("U1").Select()
